$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "TEST1"
$ws.Range("A3").Value = "TEST2"
$ws.Range("A4").Value = "TEST3"
$ws.Range("A5").Value = "TEST4"

[void]$ws.Range("G9").Select()
